$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.235.45"
$ws.Range("E2").Value = "  -0.82%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.661.75"
$ws.Range("E3").Value = "  -0.87%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.37%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'218.70"
$ws.Range("E5").Value = "  +0.99%  "

# Row 6 - XRP
$ws.Range("D6").Value = "'0.5223"
$ws.Range("E6").Value = "  -1.64%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.2670"
$ws.Range("E8").Value = "  -0.11%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.06326"
$ws.Range("E9").Value = "  -1.09%  "

# Row 10 - Solana
$ws.Range("D10").Value = "'21.08"
$ws.Range("E10").Value = "  -2.24%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.07713"
$ws.Range("E11").Value = "  -1.25%  "

# Row 12 & 13 - Coins swapped: WrappedEther/Polkadot rows exchange places
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'4.428"
$ws.Range("E12").Value = "  -1.65%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.631.53"
$ws.Range("E13").Value = "  -2.73%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "1.889.45"
$ws.Range("E14").Value = "  -0.85%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "'0.5475"
$ws.Range("E15").Value = "  -1.61%  "

# Row 16 - ShibaInu
$ws.Range("D16").Value = "0.0$([char]8325)8246"
$ws.Range("E16").Value = "  -1.18%  "

# Row 17 - Litecoin
$ws.Range("D17").Value = "'64.90"
$ws.Range("E17").Value = "  -1.16%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "26.264.92"
$ws.Range("E18").Value = "  -0.79%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  +0.35%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "'4.656"
$ws.Range("E20").Value = "  -2.47%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'194.02"
$ws.Range("E21").Value = "  -0.68%  "

# Row 22 - Avalanche
$ws.Range("D22").Value = "'10.14"
$ws.Range("E22").Value = "  -2.17%  "

# Row 23 - Chainlink
$ws.Range("D23").Value = "'6.077"
$ws.Range("E23").Value = "  -4.10%  "

# Row 24 - BinanceUSD
$ws.Range("D24").Value = "'1.007"
$ws.Range("E24").Value = "  +0.50%  "

# Row 25 - Monero
$ws.Range("D25").Value = "'138.74"
$ws.Range("E25").Value = "  -3.76%  "

# Row 26 - Stellar
$ws.Range("E26").Value = "  -2.94%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "'7.229"
$ws.Range("E27").Value = "  -2.65%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "'16.11"
$ws.Range("E28").Value = "  -1.06%  "

# Row 29 - Toncoin
$ws.Range("E29").Value = "  -1.73%  "

# Row 30 - Hedera
$ws.Range("D30").Value = "'0.05995"
$ws.Range("E30").Value = "  -2.34%  "

# Row 31 - PancakeSwap
$ws.Range("D31").Value = "'1.282"
$ws.Range("E31").Value = "  +0.64%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").Value = "'3.618"
$ws.Range("E32").Value = "  -0.14%  "

# Row 33 - Filecoin
$ws.Range("D33").Value = "'3.312"
$ws.Range("E33").Value = "  -3.95%  "

# Row 34 - LidoDAOToken
$ws.Range("D34").Value = "'1.631"
$ws.Range("E34").Value = "  -3.57%  "

# Row 35 - ARBITRUM
$ws.Range("D35").Value = "'0.9799"
$ws.Range("E35").Value = "  -2.54%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  -0.39%  "

# Row 37 - MXToken
$ws.Range("D37").Value = "'2.784"
$ws.Range("E37").Value = "  +0.08%  "

# Row 38 - ImmutableX
$ws.Range("D38").Value = "'0.5887"
$ws.Range("E38").Value = "  +2.79%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  -2.83%  "

# Row 40 - FraxShare
$ws.Range("D40").Value = "'5.946"
$ws.Range("E40").Value = "  -1.59%  "

# Row 41 - TrustWalletToken
$ws.Range("E41").Value = "  +0.20%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  +0.32%  "

# Row 43 - Maker
$ws.Range("D43").Value = "1.030.11"
$ws.Range("E43").Value = "  -4.00%  "

# Row 44 - Quant
$ws.Range("D44").Value = "'99.70"
$ws.Range("E44").Value = "  -0.36%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").Value = "1.802.62"

# Row 46 - BabyDogeCoin
$ws.Range("E46").Value = "  +5.69%  "

# Row 47 - Aave
$ws.Range("D47").Value = "'57.20"
$ws.Range("E47").Value = "  +0.19%  "

# Row 48 - Frax
$ws.Range("E48").Value = "  +0.23%  "

# Row 49 - EnergySwap
$ws.Range("D49").Value = "'8.074"
$ws.Range("E49").Value = "  -0.87%  "

# Row 50 - Cronos
$ws.Range("D50").Value = "'0.05188"
$ws.Range("E50").Value = "  -0.32%  "

# Row 51 - RenderToken
$ws.Range("E51").Value = "  -0.21%  "
